$d = $word.ActiveDocument

# --- Paragraph 1 originally has two runs split by the "_GoBack" bookmark:
#   run1: "I wrote the initial code to run on the PC which acted as a proof
#          of concept for the algorithm we would use. I also ported the
#          code to the PS3, writing the"                     (154 chars)
#   <bookmarkStart/bookmarkEnd _GoBack>
#   run2: " serial code to run on multiple SPU's. I also built the basic
#          build system using a basic python script which would allow us
#          to toggle what we built as sometimes we only wanted to build
#          the SPU programs not the PPU executable. Finally, I wrote much
#          of the common classes used to deal with files logging and high
#          resolutions timers."                               (331 chars)
#
# Total paragraph length = 154 + 331 = 485

$run1Len = 154
$paraLen = 485

# --- Step 1: rewrite run2 (after the bookmark) in place, keeping it a
# single run with no leading/trailing whitespace.
$tail = $d.Range($run1Len, $paraLen)
$tail.Text = "will certainly help me in future."

# --- Step 2: clear out the old run1 text entirely (this removes the run).
$head = $d.Range(0, $run1Len)
$head.Text = ""

# --- Step 3: rebuild run1 as many separate runs by repeatedly inserting
# at the very start of the document (position 0). Inserting there always
# creates a brand-new run rather than merging into a neighbour, so we
# insert the segments in reverse order to end up with the correct
# left-to-right order.
$segments = @(
    "During this project ",
    "I ",
    "did a lot. I wrote the initial",
    " PC",
    " code and then ported the code to the PS3",
    " to run on multiple SPU’s. I also built the basic build sys",
    "tem using a basic python script. I also ",
    "wrote much of the common classes used to deal with files",
    ",",
    " logging and high resolutions timers.",
    " I’ve been equipped with an enhanced knowledge of the PS3 but this project has also made me better a",
    "t",
    " programming for unique hardware offerings. Also the low level nature of the project has also helped me brush up on my low level programming skills",
    ". These skills",
    " "
)

for ($i = $segments.Length - 1; $i -ge 0; $i--) {
    $r = $d.Range(0, 0)
    $r.InsertBefore($segments[$i])
}

Write-Output $d.Paragraphs.Item(1).Range.Text
